$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table in B2:J16 was re-expressed in different units (values divided
# by 1000, i.e. multiplied by 1e-3). Column A (row labels) and row 1 (column
# headers) are left untouched.
$dataRange = $ws.Range("B2:J16")
foreach ($cell in $dataRange.Cells) {
    $cell.Value2 = $cell.Value2 * 0.001
}

# Update the active selection to match the saved view state (I24 instead of I25).
$ws.Range("I24").Select()
